# Day 111 to 115 Binary Tree Completed
# Marks rows 193-211 (column C, the "Done [yes or no]" status column) as
# completed ("yes"), except row 197 which is marked "NEED BST" (a follow-up
# reminder), matching the various alternating fill-color styles already
# used elsewhere in the sheet for the "yes" marker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells already carrying each fill-colour style used by the
# "yes" marker elsewhere in the sheet, so copying reuses the existing
# style definitions instead of creating new ones.
$refFill3 = $ws.Range("C13")   # fillId 3 -> style index 12
$refFill4 = $ws.Range("C16")   # fillId 4 -> style index 13
$refFill5 = $ws.Range("C20")   # fillId 5 -> style index 14
$refFill6 = $ws.Range("C22")   # fillId 6 -> style index 15
$refFill10 = $ws.Range("C240") # fillId 10 -> style index 20 ("NEED ..." marker)

# row -> reference style cell to copy (value + format) onto column C
$rowStyleRef = @{
    193 = $refFill4
    194 = $refFill3
    195 = $refFill6
    196 = $refFill5
    198 = $refFill3
    199 = $refFill3
    200 = $refFill5
    201 = $refFill3
    202 = $refFill6
    203 = $refFill6
    204 = $refFill6
    205 = $refFill3
    206 = $refFill5
    207 = $refFill6
    208 = $refFill3
    209 = $refFill3
    210 = $refFill3
    211 = $refFill6
}

foreach ($row in $rowStyleRef.Keys) {
    $src = $rowStyleRef[$row]
    $dst = $ws.Cells.Item($row, 3)
    $src.Copy($dst)
    $dst.Value2 = "yes"
}

# Row 197 gets the "NEED BST" marker instead, using the same fill style
# as the other "NEED ..." reminders in the sheet (e.g. "NEED TREE").
$dst197 = $ws.Cells.Item(197, 3)
$refFill10.Copy($dst197)
$dst197.Value2 = "NEED BST"

# Update the sheet view: scrolled position and the current selection.
# (Activating/selecting a single cell inside the range would collapse the
# selection back down to that cell, so only the range Select is issued;
# the active cell resolves to the range's anchor corner.)
$win = $excel.ActiveWindow
$win.ScrollRow = 196
$win.ScrollColumn = 1

$ws.Range("A177:C211").Select()
